$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.800000011920929
$ws.Range("F2").Value = 0.7499999999999999
$ws.Range("J2").Value = 0.6

# Row 3
$ws.Range("D3").Value = 0.9629999995231628
$ws.Range("F3").Value = 0.9615784008307372
$ws.Range("J3").Value = 0.926

# Row 4
$ws.Range("D4").Value = 0.9819999933242798
$ws.Range("F4").Value = 0.9816700610997963
$ws.Range("J4").Value = 0.964

# Row 5
$ws.Range("D5").Value = 0.9825000166893005
$ws.Range("F5").Value = 0.9821882951653944
$ws.Range("J5").Value = 0.965

# Row 6
$ws.Range("D6").Value = 0.9555000066757202
$ws.Range("F6").Value = 0.9573550551030187
$ws.Range("H6").Value = 0.9190432382704692
$ws.Range("L6").Value = 0.08799999999999999

# Row 7
$ws.Range("F7").Value = 0.9994997498749374
$ws.Range("H7").Value = 1
$ws.Range("J7").Value = 0.999
$ws.Range("L7").Value = 0

# Row 8
$ws.Range("D8").Value = 0.9994999766349792
$ws.Range("F8").Value = 0.9994997498749374
$ws.Range("H8").Value = 1
$ws.Range("L8").Value = 0

# Row 9
$ws.Range("D9").Value = 0.9994999766349792
$ws.Range("F9").Value = 0.9995002498750626
$ws.Range("H9").Value = 0.999000999000999
$ws.Range("L9").Value = 0.001
